$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.409.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "'2.435.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.83%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'565.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("D6").Value = "'145.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.82%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'0.534"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.62%  "

$ws.Range("D9").Value = "'0.111"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.49%  "

$ws.Range("D10").Value = "'0.155"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.46%  "

$ws.Range("E11").Value = "  -1.38%  "

$ws.Range("D12").Value = "'0.351"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "'26.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.18%  "

$ws.Range("D14").Value = "'0.0000179"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.87%  "

$ws.Range("D16").Value = "'62.243.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("D17").Value = "'2.439.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").Value = "'324.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").Value = "'6.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.72%  "

$ws.Range("D21").Value = "'4.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").Value = "'67.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.48%  "

$ws.Range("D24").Value = "'1.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.70%  "

$ws.Range("D25").Value = "'8.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.94%  "

$ws.Range("D26").Value = "'565.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.07%  "

$ws.Range("D27").Value = "'0.0₃0975"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.21%  "

$ws.Range("D28").Value = "'2.553.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("D30").Value = "'8.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.70%  "

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").Value = "'0.149"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D34").Value = "'1.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("D35").Value = "'4.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.23%  "

$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("E38").Value = "  +0.43%  "

$ws.Range("D39").Value = "'18.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("D40").Value = "'149.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.92%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("D43").Value = "'2.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.79%  "

$ws.Range("D44").Value = "'148.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("E45").Value = "  +1.37%  "

$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("D47").Value = "'20.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.60%  "

$ws.Range("D48").Value = "'0.600"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("D49").Value = "'0.0928"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.21%  "

$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("D51").Value = "'11.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.76%  "
